# Add a row 12 that aggregates each metric column (A:AH) over rows 2:11 as
# "ROUND(AVERAGE,3) ± ROUND(STDEV.P,3)", matching the commit's "Add files via
# upload" refresh of bigru_result_full_ce.xlsx, then restore the author's
# on-screen view (scrolled right to column L, zoomed to 85%, cell AB4 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(
    "A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z",
    "AA","AB","AC","AD","AE","AF","AG","AH"
)

foreach ($col in $cols) {
    $cell = "$col" + "12"
    $rangeTop = "$col" + "2"
    $rangeBottom = "$col" + "11"
    $formula = "=ROUND(AVERAGE($rangeTop`:$rangeBottom),3) &" + '"' + [char]0x00B1 + '"' + "& ROUND(_xlfn.STDEV.P($rangeTop`:$rangeBottom),3)"
    $ws.Range($cell).Formula = $formula
}

# Restore the saved view state: scrolled so column L is left-most visible,
# zoomed to 85%, with AB4 the active/selected cell.
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.Zoom = 85
$ws.Range("AB4").Select()
